# Fruta / hortaliza, semanal
# Insert a new weekly record at row 240, pushing the existing rows 240-287
# down to 241-288 (dimension grows from A1:R287 to A1:R288).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 240:287 down by inserting a fresh row above the current row 240.
$ws.Rows("240:240").Insert()

# Populate the newly inserted row 240 with the new observation.
$ws.Cells.Item(240, 1).Value = 4
$ws.Cells.Item(240, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(240, 3).Value = "Los Lagos"
$ws.Cells.Item(240, 4).Value = 44637
$ws.Cells.Item(240, 5).Value = 10
$ws.Cells.Item(240, 6).Value = 100112008
$ws.Cells.Item(240, 7).Value = "Coliflor"
$ws.Cells.Item(240, 8).Value = "Sin especificar"
$ws.Cells.Item(240, 9).Value = "Primera"
$ws.Cells.Item(240, 10).Value = 250
$ws.Cells.Item(240, 11).Value = 1600
$ws.Cells.Item(240, 12).Value = 1600
$ws.Cells.Item(240, 13).Value = 1600
$ws.Cells.Item(240, 14).Value = "`$/unidad"
$ws.Cells.Item(240, 15).Value = "Región Metropolitana"
$ws.Cells.Item(240, 16).Value = 1600
$ws.Cells.Item(240, 17).Value = 1
$ws.Cells.Item(240, 18).Value = "Hortaliza"
